$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update frequency counts in column C (rows 2-10)
$ws.Range("C2").Value = 2995
$ws.Range("C3").Value = 2918
$ws.Range("C4").Value = 2159
$ws.Range("C5").Value = 1288
$ws.Range("C6").Value = 1186
$ws.Range("C7").Value = 677
$ws.Range("C8").Value = 613
$ws.Range("C9").Value = 451
$ws.Range("C10").Value = 426

# Row 11: category combination changed from "Textiles & Cozy Items" / "Textiles & Cozy Items"
# to "Seasonal & Holidays" / "Home Decor", with updated frequency
$ws.Range("A11").Value = "Seasonal & Holidays"
$ws.Range("B11").Value = "Home Decor"
$ws.Range("C11").Value = 397
